# Split the run "frousse " into "fro" + "n" + "sse " where the
# inserted "n" carries minimal (default) run formatting (just rtl=0),
# matching the pattern used elsewhere in this document for inline
# corrections/insertions.

$d = $word.ActiveDocument

# First, locate the "frousse " run and replace its text with the
# expanded form "fronsse " using Find/Replace (this keeps everything
# in one run with the original formatting).
$range = $d.Content
$range.Find.Execute("frousse ", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "fronsse ", 2)

# Now find the newly-inserted "n" (the one between "fro" and "sse ")
# and give it minimal formatting by clearing the character formatting
# that differentiates it from the Normal style default.
$search = $d.Content
$search.Find.ClearFormatting()
$search.Find.Execute("fronsse ", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)

if ($search.Find.Found) {
    # $search now covers "fronsse "; narrow down to just the "n" at
    # position 3 (0-based) within that range.
    $start = $search.Start + 3
    $end = $start + 1
    $nRange = $d.Range($start, $end)
    $nRange.Font.Reset()
}
